# Insert a new data row at spreadsheet row 116 (this pushes the existing
# rows 116:221 down to 117:222, matching the rest of the sheet's weekly
# cadence) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 116:221 down by one row.
$ws.Rows.Item(116).Insert()

# Fill in the newly inserted row 116 with this week's values.
$ws.Range("A116").Value = 3
$ws.Range("B116").Value = "Femacal de La Calera"
$ws.Range("C116").Value = "Coquimbo"
$ws.Range("D116").Value = 44484
$ws.Range("E116").Value = 5
$ws.Range("F116").Value = 100112040
$ws.Range("G116").Value = "Cilantro"
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 180
$ws.Range("K116").Value = 2300
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = 2400
$ws.Range("N116").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O116").Value = "Provincia de Quillota"
$ws.Range("P116").Value = 800
$ws.Range("Q116").Value = 3
$ws.Range("R116").Value = "Hortaliza"

# Match the date-time number format used by the rest of column D.
$ws.Range("D116").NumberFormat = $ws.Range("D115").NumberFormat
